$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D" = 0.0394
    "E" = 0.114
    "K" = 13.5
    "L" = 0.08794788273615635
    "M" = 0.399
    "N" = 0.005683760683760684
    "O" = 0.02955555555555556
    "P" = 0.399
    "Q" = 0.005683760683760684
    "R" = 0.02955555555555556
    "U" = 827.2
    "V" = 11.78347578347578
    "W" = 0.07605633802816901
    "X" = 0.2238681304565773
    "Y" = -0.1478117924284083
    "Z" = 5.185810810810807
    "AB" = 0.04171459907877991
    "AC" = -0.04171459907877991
    "AD" = 577.5
    "AF" = 577.5
    "AG" = -249.7
    "AH" = 0.8916164891153311
    "AI" = 0.7120838471023427
    "AJ" = 1.39108635097493
    "AK" = 15.41358024691354
}

foreach ($row in 2..3) {
    foreach ($col in $updates.Keys) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
